{"js": "// Insert a new row \"Code_Effet_a_obtenir\" / \"CISU-Code_Effet_a_obtenir-v24.06.19\" / \"CISU\"\n// right after the \"ISO 3166\" (ISO3166-2) row in the summary table.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Locate the row whose first cell equals \"ISO3166-2\" (the \"ISO 3166\" row).\nconst items = rows.items;\nfor (const row of items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nlet targetRow = null;\nfor (const row of items) {\n  const cells = row.cells.items;\n  if (cells.length > 0 && cells[0].value === \"ISO3166-2\") {\n    targetRow = row;\n    break;\n  }\n}\n\nif (!targetRow) {\n  throw new Error(\"Could not find the ISO3166-2 row to insert after.\");\n}\n\ntargetRow.insertRows(\"After\", 1, [\n  [\"Code_Effet_a_obtenir\", \"CISU-Code_Effet_a_obtenir-v24.06.19\", \"CISU\"]\n]);\n\nawait context.sync();\n", "ps1": "# Insert a new row \"Code_Effet_a_obtenir\" / \"CISU-Code_Effet_a_obtenir-v24.06.19\" / \"CISU\"\n# right after the \"ISO 3166\" (ISO3166-2) row in the summary table.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"ISO3166-2\")\nif (-not $found) {\n    throw \"Could not find the ISO3166-2 row to insert after.\"\n}\n\n$cell = $rng.Cells.Item(1)\n$table = $cell.Tables.Item(1)\n$rowIdx = $cell.RowIndex\n\nif ($rowIdx -lt $table.Rows.Count) {\n    $beforeRow = $table.Rows.Item($rowIdx + 1)\n    $newRow = $table.Rows.Add($beforeRow)\n} else {\n    $newRow = $table.Rows.Add()\n}\n\n$newRow.Cells.Item(1).Range.Text = \"Code_Effet_a_obtenir\"\n$newRow.Cells.Item(2).Range.Text = \"CISU-Code_Effet_a_obtenir-v24.06.19\"\n$newRow.Cells.Item(3).Range.Text = \"CISU\"\n"}
